$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A26").Value = "SUTD"
$ws.Range("B26").Value = "103.963175,1.340793"

$ws.Range("B26").Select()
